# Update the "F" column (want-to-go counts) across all four sheets of the
# workbook, matching the regenerated data snapshot (gh-pages output at
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1231
$ws.Cells.Item(4, 6).Value = 50
$ws.Cells.Item(5, 6).Value = 1808
$ws.Cells.Item(6, 6).Value = 1739
$ws.Cells.Item(7, 6).Value = 6272
$ws.Cells.Item(8, 6).Value = 132
$ws.Cells.Item(9, 6).Value = 1867
$ws.Cells.Item(10, 6).Value = 495
$ws.Cells.Item(12, 6).Value = 25
$ws.Cells.Item(15, 6).Value = 44
$ws.Cells.Item(16, 6).Value = 7268
$ws.Cells.Item(19, 6).Value = 178
$ws.Cells.Item(21, 6).Value = 1724
$ws.Cells.Item(24, 6).Value = 23
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(28, 6).Value = 1659
$ws.Cells.Item(29, 6).Value = 780
$ws.Cells.Item(30, 6).Value = 341
$ws.Cells.Item(33, 6).Value = 68
$ws.Cells.Item(35, 6).Value = 3904

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(13, 6).Value = 4
$ws.Cells.Item(14, 6).Value = 25
$ws.Cells.Item(20, 6).Value = 30
$ws.Cells.Item(21, 6).Value = 62

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 2266

# Sheet "全部类型" (All types - combined listing)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 2266
$ws.Cells.Item(5, 6).Value = 1231
$ws.Cells.Item(7, 6).Value = 50
$ws.Cells.Item(10, 6).Value = 1818
$ws.Cells.Item(12, 6).Value = 1739
$ws.Cells.Item(13, 6).Value = 6272
$ws.Cells.Item(14, 6).Value = 132
$ws.Cells.Item(15, 6).Value = 1867
$ws.Cells.Item(18, 6).Value = 495
$ws.Cells.Item(20, 6).Value = 25
$ws.Cells.Item(23, 6).Value = 44
$ws.Cells.Item(24, 6).Value = 7269
$ws.Cells.Item(27, 6).Value = 178
$ws.Cells.Item(29, 6).Value = 1724
$ws.Cells.Item(31, 6).Value = 23
$ws.Cells.Item(34, 6).Value = 1659
$ws.Cells.Item(35, 6).Value = 4
$ws.Cells.Item(36, 6).Value = 341
$ws.Cells.Item(38, 6).Value = 25
$ws.Cells.Item(43, 6).Value = 3904
$ws.Cells.Item(44, 6).Value = 30
$ws.Cells.Item(45, 6).Value = 62
